$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are treated as text so formatted
# numeric-looking strings (e.g. "1.00", "0.0490") are not coerced
# into numbers by Excel, losing their exact display text.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.351.78'
$ws.Range("E2").Value = '  +5.69%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.540.19'
$ws.Range("E3").Value = '  +2.03%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.50%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '416.63'
$ws.Range("E5").Value = '  +0.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.48'
$ws.Range("E6").Value = '  -1.64%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.647'
$ws.Range("E7").Value = '  +3.37%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.533.14'
$ws.Range("E8").Value = '  +1.84%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.998'
$ws.Range("E9").Value = '  -0.13%  '
$ws.Range("E10").Value = '  +6.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.177'
$ws.Range("E11").Value = '  +25.25%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000333'
$ws.Range("E12").Value = '  +51.92%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.41'
$ws.Range("E13").Value = '  -1.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.91'
$ws.Range("E14").Value = '  +2.04%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.103.98'
$ws.Range("E15").Value = '  +1.85%  '
$ws.Range("E16").Value = '  -0.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '20.10'
$ws.Range("E17").Value = '  -2.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.502.55'
$ws.Range("E18").Value = '  +0.93%  '
$ws.Range("E19").Value = '  +3.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.38'
$ws.Range("E20").Value = '  -3.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '66.230.59'
$ws.Range("E21").Value = '  +5.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '456.24'
$ws.Range("E22").Value = '  -3.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '89.12'
$ws.Range("E23").Value = '  -2.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.17'
$ws.Range("E24").Value = '  -3.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.98'
$ws.Range("E25").Value = '  -3.78%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.34'
$ws.Range("E26").Value = '  +0.45%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.88'
$ws.Range("E27").Value = '  -6.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '34.03'
$ws.Range("E28").Value = '  +1.75%  '
$ws.Range("E29").Value = '  +0.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.77'
$ws.Range("E30").Value = '  +4.99%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.33'
$ws.Range("E31").Value = '  +2.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.116'
$ws.Range("E32").Value = '  +3.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.21'
$ws.Range("E33").Value = '  -5.63%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.158'
$ws.Range("E34").Value = '  -5.21%  '
$ws.Range("E35").Value = '  -0.19%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '38.73'
$ws.Range("E36").Value = '  -5.79%  '
$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '56.45'
$ws.Range("E37").Value = '  -2.27%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0490'
$ws.Range("E38").Value = '  -0.70%  '
$ws.Range("B39").Value = 'PEPE'
$ws.Range("C39").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0₃0772'
$ws.Range("E39").Value = '  +36.26%  '
$ws.Range("E40").Value = '  +8.85%  '
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.71'
$ws.Range("E42").Value = '  +0.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.95'
$ws.Range("E43").Value = '  -3.71%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '146.53'
$ws.Range("E44").Value = '  +0.89%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.33'
$ws.Range("E45").Value = '  -2.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.22'
$ws.Range("E46").Value = '  -4.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.306'
$ws.Range("E47").Value = '  -5.78%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.96'
$ws.Range("E48").Value = '  -5.97%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.30'
$ws.Range("E49").Value = '  -4.50%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.59'
$ws.Range("E50").Value = '  +10.24%  '
$ws.Range("B51").Value = 'Celestia'
$ws.Range("C51").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '15.35'
$ws.Range("E51").Value = '  -6.81%  '
